$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the biography paragraph
#    to the very start of the document (right before "Author Biography").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
# A zero-length range placed at absolute position 0 is mishandled by the
# bookmark engine (it silently grows to cover the first run), so nudge
# in from a throw-away character: insert it at position 0, add the
# bookmark right after it (now a perfectly interior, collapsed range),
# then delete the throw-away character again. The bookmark stays
# collapsed at position 0 once its neighbour is removed.
$throwaway = $d.Range(0, 0)
$throwaway.InsertBefore("X")
$d.Bookmarks.Add("_GoBack", $d.Range(1, 1))
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------
# 2. Edit the biography paragraph text.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Public and Social Policy at Saint Louis University", $true, $false, $false, $false, $false, $true, 1, $false, "Public and Social Policy (PSP) at Saint Louis University", 2)

$find.ClearFormatting()
$find.Execute("He earned an MBA with concentrations", $true, $false, $false, $false, $false, $true, 1, $false, "He earned a master of business administration (MBA) with concentrations", 2)

$find.ClearFormatting()
$find.Execute("technology transfer, management consulting, product management, and project management.", $true, $false, $false, $false, $false, $true, 1, $false, "technology transfer, product management, management consulting, and project management.", 2)

# ---------------------------------------------------------------------
# 3. The section now starts counting pages at 1 instead of 2.
# ---------------------------------------------------------------------
$section = $d.Sections(1)
$section.Headers(1).PageNumbers.StartingNumber = 1

# ---------------------------------------------------------------------
# 4. The header's PAGE field cached display text must match (it shows
#    the now-stale "2" cached the last time the field was computed).
#    Target the field result specifically via Find scoped to the header
#    range so the field-code runs aren't disturbed.
# ---------------------------------------------------------------------
$headerRange = $section.Headers(1).Range
$headerFind = $headerRange.Find
$headerFind.ClearFormatting()
$headerFind.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)
